$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$orgs = @(
    @("Akershus universitetssykehus HF", 29),
    @("Diakonhjemmet Sykehus", 12),
    @("Finnmarkssykehuset HF", 2),
    @("Haraldsplass Diakonale Sykehus", 3),
    @("Helgelandssykehuset HF", 0),
    @("Helse Bergen HF", 10),
    @("Helse Fonna HF", 4),
    @("Helse Førde HF", 2),
    @("Helse Møre og Romsdal", 8),
    @("Helse Nord-Trøndelag", 4),
    @("Helse Stavanger HF", 5),
    @("Lovisenberg Diakonale Sykehus", 12),
    @("Nordlandssykehuset HF", 5),
    @("Oslo universitetssykehus HF", 41),
    @("Sørlandet sykehus HF", 8),
    @("St. Olavs hospital", 9),
    @("Sunnaas Sykehus HF", 0),
    @("Sykehuset Innlandet HF", 10),
    @("Sykehuset i Vestfold HF", 6),
    @("Sykehuset Østfold HF", 14),
    @("Sykehuset Telemark HF", 3),
    @("Universitetssykehuset Nord-Norge HF", 6),
    @("Vestre Viken HF", 21)
)

$startRow = 807
$dateSerial = 43933
# Convert Excel serial date (1899-12-30 epoch) to a real date for the Value assignment.
$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0
$rowDate = $epoch.AddDays($dateSerial)

for ($i = 0; $i -lt $orgs.Length; $i++) {
    $r = $startRow + $i
    $name = $orgs[$i][0]
    $count = $orgs[$i][1]

    $ws.Cells.Item($r, 1).Value = $rowDate
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $count
}
